$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style (bold, centered, bordered) from an existing header cell (E1) to F1:H1
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Boolean flag values for Outliers_MAD columns (KNN, SVM, RF) per row
$values = @{
    2  = @($false, $false, $false)
    3  = @($false, $false, $false)
    4  = @($false, $false, $false)
    5  = @($false, $false, $false)
    6  = @($false, $false, $false)
    7  = @($false, $false, $false)
    8  = @($true,  $false, $false)
    9  = @($false, $false, $false)
    10 = @($false, $false, $false)
    11 = @($false, $false, $false)
    12 = @($false, $false, $false)
    13 = @($true,  $true,  $false)
    14 = @($false, $false, $false)
    15 = @($false, $false, $false)
    16 = @($false, $false, $false)
    17 = @($false, $false, $false)
    18 = @($false, $false, $false)
    19 = @($false, $false, $false)
    20 = @($false, $false, $false)
    21 = @($false, $false, $false)
    22 = @($true,  $true,  $false)
    23 = @($false, $false, $false)
    24 = @($false, $false, $false)
    25 = @($false, $false, $false)
}

foreach ($row in $values.Keys) {
    $vals = $values[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
    $ws.Cells.Item($row, 8).Value = $vals[2]
}
